$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) cells C1:J1 from English to Spanish labels.
$ws.Range("C1").Value = "tipo"
$ws.Range("D1").Value = "reporte"
$ws.Range("E1").Value = "nombre_tabla"
$ws.Range("F1").Value = "tabla_origen"
$ws.Range("G1").Value = "definicion"
$ws.Range("H1").Value = "columnas_agregadas_sql"
$ws.Range("I1").Value = "columnas_usadas_para_visualizacion"
$ws.Range("J1").Value = "columnas_usadas_para_powerbi"

# The longer wrapped Spanish headers now need two lines of height.
$ws.Rows.Item(1).RowHeight = 29

# Scroll the frozen pane back to the top and select C1 (matches author's
# final view state after editing the header row).
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("C1").Select() | Out-Null
